$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "jouer" row: update "avec un compte" status from "non" to "ok"
$ws.Range("B3").Value = "ok"

# "bank compte voir" row: update both status columns from "pas encore fait" to "ok"
$ws.Range("B7").Value = "ok"
$ws.Range("C7").Value = "ok"

# "notifications" row: update both status columns from "non" to "ok"
$ws.Range("B37").Value = "ok"
$ws.Range("C37").Value = "ok"

# Move the active selection to where the author was working (C25)
$ws.Range("C25").Select()
